$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Version 1." (the one with the
# _GoBack bookmark around the trailing period). We match on the visible
# text rather than assuming paragraph index 1, so the script still finds
# the right target even if earlier paragraphs exist.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "Version\s*1\.") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs(1)
}

$r = $target.Range

# Rebuild the paragraph's run layout to match the "Version 2." wording:
#  - "Version" is split into two runs ("Versi" + "on") around the
#    spell-check proofErr markers, same as the reverted revision.
#  - " 1." becomes " 2" (still inside the bookmark span) followed by a
#    brand-new run holding the final "." after the _GoBack bookmark.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Versi</w:t></w:r>
            <w:r><w:t>on</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> 2</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t>.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$countBefore = $d.Paragraphs.Count
$r.InsertXML($xml)

if ($d.Paragraphs.Count -gt $countBefore) {
    # InsertXML materialised the replacement as a brand-new paragraph ahead
    # of the (now content-less) original one; merge them back into a single
    # paragraph by deleting the paragraph mark that now separates them so
    # the original paragraph properties/rsids are preserved on the result.
    $p1End = $d.Paragraphs(1).Range.End
    $d.Range($p1End - 1, $p1End).Delete()
}
